$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.963.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.911.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4588"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3822"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07723"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9811"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.900.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.943"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.664"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07044"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "83.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009462"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.74%  "
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.941.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.319"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.095"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.657"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.850"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09295"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8689"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.079"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.99%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.248"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.28%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.028"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05726"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.155"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02041"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5502"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.88%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.400"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.73%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1753"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.841"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.73%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.317"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5186"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.22%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06867"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("B47").Value = "PEPE"
$ws.Range("C47").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000002619"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.80%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.058"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.83%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.782"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2860"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.42%  "
